$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the two paragraphs that currently hold:
#   1) the bookmarked "git commit -m """ paragraph
#   2) the following "git push origin main" paragraph
# "push origin main" is unique in the document, so we anchor on it
# and take the paragraph immediately before it as the bookmark one.
# ------------------------------------------------------------------
$finder = $d.Content
$finder.Find.ClearFormatting()
$ok = $finder.Find.Execute("push origin main", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not locate 'push origin main' paragraph"
}

$pushPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($finder.Start -ge $cand.Range.Start -and $finder.Start -lt $cand.Range.End) {
        $pushPara = $cand
        $pushIndex = $i
        break
    }
}
if ($pushPara -eq $null) {
    throw "Could not map found text back to a paragraph"
}

$commitPara = $d.Paragraphs.Item($pushIndex - 1)

# Sanity check: the commit paragraph should contain the _GoBack bookmark
$bm = $d.Bookmarks.Item("_GoBack")
if (-not ($bm.Start -ge $commitPara.Range.Start -and $bm.Start -le $commitPara.Range.End)) {
    throw "Unexpected document layout: bookmark not in expected paragraph"
}

# ------------------------------------------------------------------
# Remove both paragraphs (text, runs, proofErr marks and bookmark)
# completely, including their paragraph marks.
# ------------------------------------------------------------------
$insertPos = $commitPara.Range.Start
$killRange = $d.Range($commitPara.Range.Start, $pushPara.Range.End)
$killRange.Delete()

# ------------------------------------------------------------------
# Insert the replacement paragraphs as raw WordprocessingML so that
# run / proofErr boundaries match exactly:
#   - "git commit -m """          (unchanged content, no bookmark)
#   - "git push origin main"      (unchanged content)
#   - empty paragraph
#   - "highlight and ctlr slash to comment"
#   - "rgb and colors"            (now carries the _GoBack bookmark
#                                   at the end of the paragraph)
# ------------------------------------------------------------------
$dash = [char]0x2013
$lq = [char]0x201C
$rq = [char]0x201D

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>git</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> commit ' + $dash + 'm ' + $lq + $rq + '</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>git</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> push origin main</w:t></w:r>' +
  '</w:p>' +
  '<w:p/>' +
  '<w:p>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>highlight</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ctlr</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> slash to comment</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>rgb</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and colors</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$target = $d.Range($insertPos, $insertPos)
$null = $target.InsertXML($xml)
